$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("esercizio base")

# Row totals and averages for gennaio..maggio (rows 3-7)
$ws.Range("G3").Formula = "=C3+D3+E3"
$ws.Range("H3").Formula = "=(C3+D3+E3)/3"

$ws.Range("G4:G7").Formula = "=C4+D4+E4"
$ws.Range("H4:H7").Formula = "=(C4+D4+E4)/3"

# Column totals (row 16) and averages (row 17)
$ws.Range("C16").Formula = "=SUM(C3:C14)"
$ws.Range("D16:E16").Formula = "=SUM(D3:D14)"

$ws.Range("C17").Formula = "=AVERAGE(C3:C14)"
$ws.Range("D17:E17").Formula = "=AVERAGE(D3:D14)"

# Update selection to match final state
$ws.Range("K13").Select()
